# Latest status of action items.
# Update the Disposition column (E) for the caArray/caIntegrator upgrade &
# permissions-wireframe action items, which have since completed, and
# normalize their row shading to match the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item 18 - "Upgrade Curation tier to caArray 2.4.1." -> Complete
$ws.Range("E19").Value = "Complete"

# Item 19 - "Create wireframes to depict how permissions will work..." -> Complete
$ws.Range("E20").Value = "Complete"

# Item 20 - "Work with Eve to schedule a meeting with UCSF." -> Complete
$ws.Range("E21").Value = "Complete"

# Item 21 - "Links to the caArray User's Guide..." -> Complete
$ws.Range("E22").Value = "Complete"

# Item 22 - "Create sample experiments in caArray to represent TRANSCEND use cases." -> Complete
$ws.Range("E23").Value = "Complete"

# These rows (17-22, i.e. spreadsheet rows 18-23) now match the disposition of
# the rest of the completed items above them, so bring their shading in line
# with the rest of the table's banding.
$ws.Range("A18:E23").Interior.ColorIndex = 22

# Reflect the cell the author had selected when the sheet was last saved.
$ws.Range("C28").Select()
